$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.21, 4693.743604618, 34.74120455, 4659.002400068, 25.32429377533333, 1764.728998344, 21.748509934, 1742.98048841, 26.59494055133333),
    @(0.22, 4873.862068222, 34.549828712, 4839.31223951, 25.45483019533333, 1814.294152412, 21.730126024, 1792.564026386, 26.73886206266667)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 22 + $i
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}
